# Estado de Cuenta update:
#  - Elimina antiguos EC y agrega nuevos (nuevo periodo 2509, valores actualizados)
#  - Modifica Antigua BD (VALOR MORA, Cant. Periodos)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update VALOR MORA (E11) and Cant. Periodos (F13)
$ws.Range("E11").Value = 624000
$ws.Range("F13").Value = 12

# 2) Insert a new blank row right after the last data row (26), before row 31,
#    which also shifts the signature rows (31,32) down to (32,33).
$ws.Rows.Item(27).Insert()

# 3) Populate new row 27 as a copy of row 26 (same worker, new period), then
#    update its period value to the new one (2509).
$ws.Range("B26:J26").Copy($ws.Range("B27:J27"))
$ws.Range("E27").Value = "2509"

# 4) Row 26 no longer is the last row of the table, so it should take the
#    regular "middle" row formatting (same as rows 16-25) instead of the
#    table-closing border; copy formatting only from row 25.
$ws.Range("B25:J25").Copy()
$ws.Range("B26:J26").PasteSpecial(-4122)
$excel.CutCopyMode = $false
